$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 9) to the Item table:
# ID, ItemType, ItemSubType, Level, ShowName, Desc, Icon, CoolDownTime, OverlayCount, BuyPrice, SalePrice
$ws.Cells.Item(9, 1).Value = "Equip_Weapon_1"
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(9, 5).Value = "开山斧"
$ws.Cells.Item(9, 6).Value = "开山斧武器"

# Icon is stored as text (e.g. "1017", "1018" ...), so force text format before assigning.
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = "50004"

$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 10000
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 100

# Grow the XML table ("表1") to include the new row -> updates table ref & autoFilter.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K9"))

# Update the last selection recorded in the sheet view.
$ws.Range("K13").Select()
